# NATMI re-ran its ligand-receptor scoring with an updated TPM matrix.
# Only the "Receptor average expression value" column (M) is raw input;
# every other touched column (N, O, P, Q, R, S, T) is a value NATMI
# derives from M (N = M * receptor-expressing-cell count; Q = ligand-avg *
# M; R = ligand-total * N; O/P/S/T = per-row share of the M/N column
# total). We write the resulting literal values directly (matching what
# the regenerated NATMI output file contains) so the numbers land exactly
# as the refreshed run produced them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ col letter = new value }
$updates = @{
    2 = @{ M = 0.6201396666666666; N = 1.860419;       O = 0.1891028895630684; P = 0.1891028895630684;
           Q = 0.07186509198488888; R = 0.6467858278640001; S = 0.1891028895630684; T = 0.1891028895630684 }
    3 = @{                                              O = 0.3106459956672665; P = 0.3106459956672664;
                                                         S = 0.3106459956672665; T = 0.3106459956672664 }
    4 = @{ M = 0.5912873333333333; N = 1.773862;       O = 0.1803047753684109; P = 0.1803047753684109;
           Q = 0.0685215297191111;  R = 0.616693767472;     S = 0.1803047753684109; T = 0.1803047753684109 }
    5 = @{ M = 0.4518883333333334; N = 1.355665;       O = 0.1377970063622857; P = 0.1377970063622857;
           Q = 0.05236723013777779; R = 0.4713050712400001; S = 0.1377970063622857; T = 0.1377970063622857 }
    6 = @{ M = 0.5973363333333334; N = 1.792009;       O = 0.1821493330389685; P = 0.1821493330389685;
           Q = 0.06922252010044445; R = 0.623002680904;     S = 0.1821493330389685; T = 0.1821493330389685 }
}

$colIndex = @{ M = 13; N = 14; O = 15; P = 16; Q = 17; R = 18; S = 19; T = 20 }

foreach ($row in $updates.Keys) {
    $rowUpdates = $updates[$row]
    foreach ($col in $rowUpdates.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $rowUpdates[$col]
    }
}
